$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Patient")

$ws.Range("A2").Value = "621730af-3caa-46a9-bea4-516276cb74e2"
$ws.Range("B2").Value = "CZTT072107"
$ws.Range("C2").Value = "Agnes"
$ws.Range("E2").Value = "Marinai"
$ws.Range("F2").Value = "Bertin"
$ws.Range("G2").Value = 25456
$ws.Range("I2").Value = "female"
